$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second "~TFM_INS" block (H2 Export Demand, ACT_BND on EXPH2*) is being
# moved out to its own file, so its data is removed from this sheet. The
# neighbouring TB_ELCC_* block in columns T:AE is left untouched.

# Row 16: C16 loses its "~TFM_INS" label but keeps its (blank) style.
$ws.Range("C16").ClearContents() | Out-Null

# Row 17: the whole header row C17:N17 (TimeSlice..Cset_CN) disappears,
# including its styling, and the row's thick-bottom/extra height (which
# Excel derived from that header's border) goes away with it.
$ws.Range("C17:N17").Clear() | Out-Null
$ws.Rows(17).AutoFit() | Out-Null

# Row 18: D18,E18,F18,H18 removed; M18 keeps its style but loses its value.
$ws.Range("D18:F18").Clear() | Out-Null
$ws.Range("H18").Clear() | Out-Null
$ws.Range("M18").ClearContents() | Out-Null

# Row 19: D19,E19,F19,H19 removed; M19 keeps its style but loses its value.
$ws.Range("D19:F19").Clear() | Out-Null
$ws.Range("H19").Clear() | Out-Null
$ws.Range("M19").ClearContents() | Out-Null

# Refresh the sheet's remembered selection/scroll position.
$ws.Range("P15").Select() | Out-Null
